$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column K ("최종점수") updates: +0.1 each
$ws.Range("K2").Value = 49.2
$ws.Range("K3").Value = 48.4
$ws.Range("K4").Value = 47.2
$ws.Range("K5").Value = 41.6
$ws.Range("K6").Value = 40.4

# Column N ("MACRO_SCORE") updates: new value for rows 2-6
$ws.Range("N2:N6").Value = 54.02451352198364
